$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 2.9
$ws.Range("I3").Value = 2.35
$ws.Range("K3").Value = 1.91
$ws.Range("L3").Value = 3.2
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 1.9
$ws.Range("R3").Value = 1.95
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.5
$ws.Range("U3").Value = 3.95
$ws.Range("V3").Value = 1.24
$ws.Range("W3").Value = 5
$ws.Range("X3").Value = 1.17
$ws.Range("Y3").Value = 1.57
$ws.Range("Z3").Value = 2.25
$ws.Range("AA3").Value = 2.1
$ws.Range("AB3").Value = 1.67
$ws.Range("AC3").Value = 8
$ws.Range("AF3").Value = 34
$ws.Range("AM3").Value = 6
$ws.Range("AQ3").Value = 23
